$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 495; this shifts existing rows 495-564 down to 496-565.
# The inserted row inherits the formatting (including the date style on column D)
# from the row that was previously at that position.
$ws.Rows.Item(495).Insert()

# Populate the new row 495 with the new record's values
$ws.Range("A495").Value = 3
$ws.Range("B495").Value = "Femacal de La Calera"
$ws.Range("C495").Value = "Coquimbo"
$ws.Range("D495").Value = 44984
$ws.Range("E495").Value = 5
$ws.Range("F495").Value = 100112017
$ws.Range("G495").Value = "Apio"
$ws.Range("H495").Value = "Americana (o)"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 200
$ws.Range("K495").Value = 9500
$ws.Range("L495").Value = 10000
$ws.Range("M495").Value = 9775
$ws.Range("N495").Value = "`$/docena de matas"
$ws.Range("O495").Value = "Provincia de Santiago"
$ws.Range("P495").Value = 1629
$ws.Range("Q495").Value = 6
$ws.Range("R495").Value = "Hortaliza"
